$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new note row: G5 gets the new remark string
$ws.Range("G5").Value = "消失：24、39、68、73、92、102"

# Update the active selection to H13 (as reflected in the saved view state)
$ws.Range("H13").Select()
